$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B2:H13")
foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = [math]::Round([double]$val)
    }
}
